$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Beta) unfolded values F2:N2 ---
$ws.Range("F2").Value = 11.43986483888852
$ws.Range("G2").Value = 11.06689720615159
$ws.Range("H2").Value = 11.82967442350957
$ws.Range("I2").Value = 1.959362889200075
$ws.Range("J2").Value = 1.933169988130725
$ws.Range("K2").Value = 1.984673215894531
$ws.Range("L2").Value = 0.1523228311846385
$ws.Range("M2").Value = 0.1502978599400667
$ws.Range("N2").Value = 0.1542958738628188

# --- Update existing row 3 (Gamma) unfolded values F3:N3 ---
$ws.Range("F3").Value = 0.09968007603381999
$ws.Range("G3").Value = 0.02915405872278833
$ws.Range("H3").Value = 0.1811423969066214
$ws.Range("I3").Value = 0.09062551516148228
$ws.Range("J3").Value = 0.02655197504282738
$ws.Range("K3").Value = 0.1643991811025323
$ws.Range("L3").Value = 0.106583413174184
$ws.Range("M3").Value = 0.03130706309780824
$ws.Range("N3").Value = 0.1932860895925299

# --- Add new row 4 (Beta + Gamma) ---
# Column A on rows 2/3 carries the bold/centered/bordered header-like
# style; reuse it for A4 via copy/paste-special (format only) so the new
# cell shares the same style entry instead of minting a new one.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 12.00687180793019
$ws.Range("D4").Value = 1.974477778970852
$ws.Range("E4").Value = 0.1537386519519979
$ws.Range("F4").Value = 11.53954491492234
$ws.Range("G4").Value = 11.09605126487438
$ws.Range("H4").Value = 12.0108168204162
$ws.Range("I4").Value = 2.049988404361557
$ws.Range("J4").Value = 1.959721963173552
$ws.Range("K4").Value = 2.149072396997063
$ws.Range("L4").Value = 0.2589062443588225
$ws.Range("M4").Value = 0.1816049230378749
$ws.Range("N4").Value = 0.3475819634553487
